$p = $ppt.ActivePresentation
$hm = $p.HandoutMaster
$hf = $hm.HeadersFooters
$dt = $hf.DateAndTime
try {
  $dt.UseFormat = 0
  Write-Output "UseFormat=0 ok"
} catch { Write-Output "ERROR: $_" }
try {
  $dt.Value = "10/22/2016"
  Write-Output "Value set ok"
} catch { Write-Output "ERROR Value: $_" }
